$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Gemma-7B-Instruct"
$ws.Range("B8").Value = "0.77 ± 0.42"
$ws.Range("C8").Value = "-0.03 ± 0.65"
$ws.Range("D8").Value = "0.35 ± 0.71"
$ws.Range("E8").Value = "0.01 ± 0.01"
$ws.Range("F8").Value = "0.12 ± 0.08"
$ws.Range("G8").Value = "0.01 ± 0.03"
$ws.Range("H8").Value = "0.11 ± 0.07"
$ws.Range("I8").Value = "0.12 ± 0.09"
$ws.Range("J8").Value = "0.78 ± 0.23"
$ws.Range("K8").Value = "0.79 ± 0.23"
$ws.Range("L8").Value = "0.79 ± 0.23"
$ws.Range("M8").Value = "0.8 ± 0.24"
$ws.Range("N8").Value = "0.91 ± 0.27"
$ws.Range("O8").Value = "0.08 ± 0.12"
$ws.Range("P8").Value = "0.46 ± 0.19"
$ws.Range("Q8").Value = "7.47 ± 1.25"
$ws.Range("R8").Value = "0.078 ± 0.00"
$ws.Range("S8").Value = "0.87 ± 0.26"
$ws.Range("T8").Value = "0.91 ± 0.28"
$ws.Range("U8").Value = "2.91 ± 1.36"
$ws.Range("V8").Value = "0.47 ± 0.45"
$ws.Range("W8").Value = "0.87 ± 0.26"
$ws.Range("X8").Value = "1.2 ± 0.42"
